$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "Start time"
$ws.Range("J1").Value = "Projected end time"

# Row 5: start time and projected end time (formula)
$ws.Range("I5").Value = 0.6479166666666667
$ws.Range("I5").NumberFormat = "h:mm AM/PM"
$ws.Range("J5").Formula = "=I5+(H5 * 1/24)"
$ws.Range("J5").NumberFormat = "h:mm AM/PM"

# Row 6: start time and projected end time (formula)
$ws.Range("I6").Value = 0.6479166666666667
$ws.Range("I6").NumberFormat = "h:mm AM/PM"
$ws.Range("J6").Formula = "=I6+(H6 * 1/24)"
$ws.Range("J6").NumberFormat = "h:mm AM/PM"

# Update the active selection to match the committed workbook state
[void]$ws.Range("J11").Select()
